# chore(results): update lottery results 2025-09-19T17:41:23Z
#
# Append the new Pick 4 draw for 2025-09-19 as row 3 of the Results sheet.
# Columns: A=Date, B=Game, C=Phase, D=Result, E=InsertedAt

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 ("2025-09-19") and C3 ("250919") look like a date / a number to Excel's
# auto-detection, so force the cell to Text first (matching the existing
# rows, which already store these "number/date-looking" values as plain
# text) and restore the Normal style afterwards so no stray number format
# sticks to the cell.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-09-19"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "Pick 4"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "250919"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "3-3-0-2"

# Timestamp string - the trailing offset (+04:00) keeps Excel from treating
# it as a recognized date/time value, so it is entered as plain text.
$ws.Range("E3").Value = "2025-09-19T21:41:23.583+04:00"
